# Update column G ("K") values for rows 2-23 in the active worksheet.
# These values were regenerated (per commit: "regen save_data to use K
# instead of Strike#, regen std/mean, calc and write s_vals") and are
# written here as plain literal values matching the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 3
    4  = 3
    5  = 1
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 3
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 2
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
